$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's results as a new row, mirroring the existing
# data layout (Day, Chase, Bryce, Zach). This is the daily 8 AM UTC
# update that adds one new row of standings to the bottom of the table.
$newRow = 66

$ws.Cells.Item($newRow, 1).Value = 46015
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = 141
$ws.Cells.Item($newRow, 3).Value = 156
$ws.Cells.Item($newRow, 4).Value = 146
